# Manual.docx edit: add new section headers / intro content before the
# "Tablets:" section, add a "Phones:" sub-section and an iPad Cellular
# note, and append a new set of section headers after the device list.

$d = $word.ActiveDocument

function Insert-ParasBefore($paragraph, [string]$bodyXml) {
    # Create a scratch empty paragraph immediately before $paragraph via
    # Range.InsertParagraphBefore(). That call leaves $paragraph's Range
    # anchored to its original *start* offset, so $paragraph itself now
    # denotes the new (empty) scratch paragraph -- "Tablets:" (etc.)
    # shifted down to become the next one.
    #
    # InsertXML at the scratch paragraph's *end* (collapse direction 0)
    # then replaces that scratch paragraph with the last paragraph of
    # $bodyXml, and materializes every earlier paragraph in $bodyXml as
    # new standalone paragraphs before it -- leaving the original target
    # paragraph text untouched.
    $paragraph.Range.InsertParagraphBefore()
    $insPoint = $paragraph.Range
    $insPoint.Collapse(0)

    $pkg = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
    $insPoint.InsertXML($pkg)
}

function Insert-ParasAfter($paragraph, [string]$bodyXml) {
    # Create a scratch empty paragraph immediately after $paragraph, then
    # InsertXML at its end -- same replace-the-scratch-paragraph trick,
    # but anchored on the following side.
    $paragraph.Range.InsertParagraphAfter()
    $scratch = $paragraph.Next()
    $insPoint = $scratch.Range
    $insPoint.Collapse(0)

    $pkg = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
    $insPoint.InsertXML($pkg)
}

function Find-ParagraphByText($doc, [string]$text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -eq ($text + "`r")) {
            return $doc.Paragraphs($i)
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Before the first paragraph ("Tablets:"), insert the new intro
#    section headers plus the DEVICES/Phones sub-section.
# ---------------------------------------------------------------------
$tabletsPara = Find-ParagraphByText $d "Tablets:"

$block1 = @"
<w:p><w:r><w:t>WARNING:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>DESIGN:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>DARK MODE:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>INSTALLING DEVICE IN PLANE:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>POWER / CHARGING:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>GPS RECEPTION:</w:t></w:r></w:p>
<w:p/>
<w:p/>
<w:p><w:r><w:t>DEVICES:</w:t></w:r></w:p>
<w:p><w:r><w:t>Phones:</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Just about any smartphone since they have real GPS. Does not need active cell service. Use </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>WiFi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to download app. </w:t></w:r></w:p>
"@

Insert-ParasBefore $tabletsPara $block1

# ---------------------------------------------------------------------
# 2) Right before "Google Nexus 9" (which directly follows "Tablets:"),
#    insert the iPad Cellular note.
# ---------------------------------------------------------------------
$nexusPara = Find-ParagraphByText $d "Google Nexus 9"

$block2 = @"
<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>iPad&#8217;s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> with Cellular (DO NOT USE </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Wifi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Only iPads)</w:t></w:r></w:p>
"@

Insert-ParasBefore $nexusPara $block2

# ---------------------------------------------------------------------
# 3) After the "Garmin GLO" paragraph, insert the closing section
#    headers (before the pre-existing trailing empty paragraph).
# ---------------------------------------------------------------------
$garminPara = Find-ParagraphByText $d "Garmin GLO"

$block3 = @"
<w:p/>
<w:p><w:r><w:t>&gt;99:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>TROUBLESHOOTING:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:lastRenderedPageBreak/><w:t>FAQ:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>FUTURE WORK:</w:t></w:r></w:p>
<w:p/>
"@

Insert-ParasAfter $garminPara $block3

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
